$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new header cells (H1:P1) ---
$newHeaders = @("grade_total", "grade_distance", "grade_visitation", "grade_encounters", `
                "NEVER", "RARELY", "SOMETIMES", "FREQUENTLY", "ALWAYS")
$startCol = 8  # H
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $startCol + $i).Value = $newHeaders[$i]
}

# --- Swap columns A and B for data rows 2-13 (A=month, B=year) ---
for ($r = 2; $r -le 13; $r++) {
    $aVal = $ws.Cells.Item($r, 1).Value2
    $bVal = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value = $bVal
    $ws.Cells.Item($r, 2).Value = $aVal
}

# --- Populate the new grade columns (H:P) for data rows 2-13 ---
# H=grade_total, I=grade_distance, J=grade_visitation, K=grade_encounters -> all 0
# L=NEVER, M=RARELY, N=SOMETIMES, O=FREQUENTLY, P=ALWAYS -> fixed values
$gradeValues = @(0, 0, 0, 0, 1.068, 1.077, 1.125, 1.227, 1.503)
for ($r = 2; $r -le 13; $r++) {
    for ($i = 0; $i -lt $gradeValues.Length; $i++) {
        $ws.Cells.Item($r, $startCol + $i).Value = $gradeValues[$i]
    }
}
